# Updates cryptos list values per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "0.556") need to be forced
# to text first, otherwise Excel COM coerces the string into a Double and the
# original formatted representation (trailing zeros, exact decimal text) is lost.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D15", "D17", "D18", "D19", "D21", "D25", "D26", "D27", "D29", "D32", "D36", "D39", "D40", "D43", "D45", "D47", "D48")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.366.09"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "1.837.62"
$ws.Range("E3").Value = "  +3.43%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "225.16"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "0.556"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "32.04"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  +4.10%  "
$ws.Range("D10").Value = "0.0728"
$ws.Range("E10").Value = "  +11.06%  "
$ws.Range("D11").Value = "0.0932"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "2.101.70"
$ws.Range("E12").Value = "  +3.43%  "
$ws.Range("D13").Value = "1.840.46"
$ws.Range("E13").Value = "  +3.61%  "
$ws.Range("E14").Value = "  +3.56%  "
$ws.Range("D15").Value = "10.80"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "34.368.80"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "4.34"
$ws.Range("E17").Value = "  +3.35%  "
$ws.Range("D18").Value = "69.75"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").Value = "252.21"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "0.0₃0799"
$ws.Range("E20").Value = "  +8.30%  "
$ws.Range("D21").Value = "11.24"
$ws.Range("E21").Value = "  +8.69%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  +2.32%  "
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").Value = "160.52"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("D26").Value = "16.70"
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("D27").Value = "7.28"
$ws.Range("E27").Value = "  +4.10%  "
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  +4.79%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("D32").Value = "1.22"
$ws.Range("E32").Value = "  +2.33%  "
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("D35").Value = "1.449.59"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").Value = "0.646"
$ws.Range("E36").Value = "  +3.86%  "
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").Value = "0.970"
$ws.Range("E39").Value = "  +9.66%  "
$ws.Range("D40").Value = "82.02"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "2.14"
$ws.Range("E43").Value = "  +4.55%  "
$ws.Range("E44").Value = "  +4.77%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "0.0502"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.996.01"
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("D47").Value = "1.05"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "106.72"
$ws.Range("E48").Value = "  +8.67%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("D51").Value = "0.0₆0125"
$ws.Range("E51").Value = "  +6.47%  "

# Restore the default (Normal) style on the forced-text cells so no stray
# number-format / style index is left behind compared to the original file.
foreach ($cell in $textCells) {
    $ws.Range($cell).Style = "Normal"
}
